# SPIFFS file system for html, latching relay logic
# Update the SmartSwitch BOM worksheet: revise several component rows and
# append two new rows (PANEL MOUNT LED + rocker switch entries).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Row 3: capacitor C1 -> C2 ---
$ws.Range("B3").Value = "Cap, Alu Elec, 1uf, 50v, Rad"
$ws.Range("C3").Value = "C2"
$ws.Range("D3").Value = "ESK105M050AC3AA"

# --- Row 4: RESET LED ---
$ws.Range("B4").Value = "LED Uni-Color Red 645nm 2-Pin T-1 3/4"
$ws.Range("D4").Value = "HLMP-D150"

# --- Row 5: SET LED ---
$ws.Range("B5").Value = "LED Uni-Color Red 645nm 2-Pin T-1 3/4"
$ws.Range("D5").Value = "HLMP-D150"

# --- Row 10: MOSFET part number tweak ---
$ws.Range("D10").Value = "DMN65D87"

# --- Row 11: resistor part number tweak (strip LRM marks) ---
$ws.Range("D11").Value = "ERJ-3EKF3300V"

# --- Row 14: was slide switch SW1, now voltage regulator U1 ---
$ws.Range("A14").Value = "TL750L05CKCS"
$ws.Range("B14").Value = "IC REG LINEAR 5V 150MA TO220-3"
$ws.Range("C14").Value = "U1"
$ws.Range("D14").Value = "TL750L05CKCS"

# --- Row 15: diode rework for latching relay logic ---
$ws.Range("A15").Value = "1N4148"
$ws.Range("B15").Value = "Diode: switching; SMD; 100V; 0.15A; 4ns; 400mW; Package: reel, tape"
$ws.Range("D15").Value = "1N4148W"

# --- Row 16: was external 5V terminal block, now LED GREEN panel mount ---
$ws.Range("A16").Value = "5111F5"
$ws.Range("B16").Value = "LED GREEN 5/32`" HOLE PANEL MOUNT"
$ws.Range("C16").Value = "PANEL MOUNT"
$ws.Range("D16").Value = "5111F5"
$ws.Range("E16").Value = 1

# --- New row 17: LED RED panel mount ---
$ws.Range("A17").Value = "5111F1"
$ws.Range("B17").Value = "LED RED 5/32`" HOLE PANEL MOUNT"
$ws.Range("C17").Value = "PANEL MOUNT"
$ws.Range("D17").Value = "5111F1"
$ws.Range("E17").Value = 1

# --- New row 18: rocker switch ---
$ws.Range("A18").Value = "RRA1GC1100"
$ws.Range("B18").Value = "SWITCH ROCKER DPDT 15A 125V"
$ws.Range("C18").Value = "PANEL MOUNT"
$ws.Range("D18").Value = "RRA1GC1100"
$ws.Range("E18").Value = 1

# Re-assert the quote-prefixed text style used throughout columns A:D (a
# plain .Value assignment can otherwise flip the cell to the no-quote-prefix
# variant of the style), and the matching style for the newly added
# quantity cells in column E.
$ws.Range("A2").Copy() | Out-Null
$ws.Range("A3:D18").PasteSpecial(-4122) | Out-Null
$ws.Range("A9").Copy() | Out-Null
$ws.Range("E16:E18").PasteSpecial(-4122) | Out-Null
$excel.CutCopyMode = 0

# --- Sheet view: zoom level and active selection ---
$ws.Activate() | Out-Null
$excel.ActiveWindow.Zoom = 160
$ws.Range("H12").Select() | Out-Null

# --- Column widths (closest achievable values under this runtime's
#     character-width rounding) ---
$ws.Columns.Item(1).ColumnWidth = 21.833333333333332
$ws.Columns.Item(2).ColumnWidth = 37

Write-Host "BOM updated"
